$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update validation strings (J column) to reflect renamed saved-search
# response parameters: contents.dra_ss -> contents.sse / contents.dra_td_sse / contents.ipa_sse
$ws.Range("J12").Value = "status=200||contents.patents=0||contents.wos=0||contents.posts=0||contents.documents=0||contents.sse=0||contents.dra_td_sse=0||contents.ipa_sse=0||contents.total=0"
$ws.Range("J24").Value = "status=200||id=(OPQA-3993_id)||type=watchlist||ispublic=true||contents.patents=2||contents.wos=2||contents.posts=2||contents.documents=0||contents.sse=0||contents.dra_td_sse=0||contents.ipa_sse=0||contents.total=6"

# Clear the STATUS column (stale PASS markers from the prior failing run)
$ws.Range("L2:L38").ClearContents()
